$wb = $excel.ActiveWorkbook

# Update sheet "展览" (first sheet): F2 84 -> 85, F3 8 -> 9
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 85
$ws1.Range("F3").Value = 9

# Update sheet "全部类型" (fourth sheet): F2 84 -> 85, F3 8 -> 9
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 85
$ws4.Range("F3").Value = 9
